# Apply the Nov-2023 performance-sheet update:
# a new order ("曹卓肺癌和癌旁组织对比分析") is inserted at row 10,
# pushing the previous rows 10-12 down to rows 11-13, and the
# weighting (K column) plus the monthly summary (rows 29-30) are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("绩效表")

# --- K-column weights that were previously blank ---
$ws.Range("K8").Value = 0.083
$ws.Range("K9").Value = 0.117

# --- Row 10: becomes the newly inserted order ---
$ws.Range("C10").Value = "曹卓肺癌和癌旁组织对比分析"
$ws.Range("E10").Value = ""
$ws.Range("G10").Value = "肺癌和癌旁组织单细胞数据对比分析"

# --- Row 11: now holds what used to be row 10's content ---
$ws.Range("C11").Value = "IN2023110603"
$ws.Range("D11").Value = "固定业务"
$ws.Range("E11").Value = "3-4"
$ws.Range("G11").Value = "OCTA 在糖尿病视网膜病变中的应用"
$ws.Range("K11").Value = 0.25

# --- Row 12: now holds what used to be row 11's content ---
$ws.Range("C12").Value = "A2023060507"
$ws.Range("D12").Value = "其他业务"
$ws.Range("E12").Value = "2-3"
$ws.Range("G12").Value = "脓毒症肠损伤联合肠道菌与代谢物分析"
$ws.Range("I12").Value = "完成"
$ws.Range("K12").Value = 0.083

# --- Row 13: brand-new row holding what used to be row 12's content ---
$ws.Range("A13").Value = "黄礼闯"
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = "曹卓补充订单"
$ws.Range("D13").Value = "其他业务"
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = "曹卓交付三个订单所需数据"
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = "完成"
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = 0.083

# --- Monthly summary rows ---
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 0.75
$ws.Range("I29").Value = "0.25+0.25+0.25=0.75"
$ws.Range("J29").Value = 1.116

$ws.Range("H30").Value = 0.366
$ws.Range("I30").Value = "0.083+0.117+0.083+0.083=0.366"
